# ------------------------------------------------------------------
# Applies the "const factor" / dual-theta-series edit to Sheet1 and
# its chart, per the target diff.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# -----------------------------------------------------------------
# 1. theta (column H) series now starts at 1 and steps by 3 (was
#    starting at 0 and stepping by 5) for each of the three blocks
#    of rows (3-8, 9-14, 15-20).
# -----------------------------------------------------------------

# Block rows 9-14 : H9 literal, H10..H14 formulas "=prev+3"
$ws.Range("H9").Value = 1
$ws.Range("H10").Formula = "=H9+3"
$ws.Range("H11").Formula = "=H10+3"
$ws.Range("H12").Formula = "=H11+3"
$ws.Range("H13").Formula = "=H12+3"
$ws.Range("H14").Formula = "=H13+3"

# Block rows 15-20 : H15 literal, H16..H20 formulas "=prev+3"
$ws.Range("H15").Value = 1
$ws.Range("H16").Formula = "=H15+3"
$ws.Range("H17").Formula = "=H16+3"
$ws.Range("H18").Formula = "=H17+3"
$ws.Range("H19").Formula = "=H18+3"
$ws.Range("H20").Formula = "=H19+3"

# -----------------------------------------------------------------
# 2. column K ("total perp_cost and angle_cost") formula changed:
#    2 * (LOG(2 * E * SIN(RADIANS(H)), 2))
#       -> 2 * (2 * LOG(E * SIN(RADIANS(H)), 2))
#    applied uniformly for every data row (3-20).
# -----------------------------------------------------------------
for ($r = 3; $r -le 20; $r++) {
    $ws.Range("K$r").Formula = "=2 * (2 * LOG(E$r * SIN(RADIANS(H$r)), 2))"
}

# -----------------------------------------------------------------
# 3. new helper cells: M5 label, M6 constant factor value.
# -----------------------------------------------------------------
$ws.Range("M5").Value = "const factor"
$ws.Range("M6").Value = 8

# -----------------------------------------------------------------
# 4. column S ("par_cost") formula now adds the constant factor in
#    M6: O+Q -> O+Q+$M$6, for every data row (3-20).
# -----------------------------------------------------------------
for ($r = 3; $r -le 20; $r++) {
    $ws.Range("S$r").Formula = "=O$r+Q$r + `$M`$6"
}

# -----------------------------------------------------------------
# 5. chart: add the two new series for the theta=10 block, and move
#    the chart object to its new anchor position.
# -----------------------------------------------------------------
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$sc = $chart.SeriesCollection()

$noParCost10 = $sc.NewSeries()
$noParCost10.Name = "no_par_cost_10"
$noParCost10.Values = "=Sheet1!`$I`$9:`$I`$14"
$noParCost10.MarkerStyle = "none"

$parCost10 = $sc.NewSeries()
$parCost10.Name = "par_cost_10"
$parCost10.Values = "=Sheet1!`$S`$9:`$S`$14"
$parCost10.MarkerStyle = "none"

# move/resize the chart object (from col B/row5 area to col L/row19 area)
$co.Left = 94.4375
$co.Top = 66.75
$co.Width = 443.5
$co.Height = 216.0

# -----------------------------------------------------------------
# 6. selection moved to M7
# -----------------------------------------------------------------
$ws.Range("M7").Select()
